$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated stat values (B:G) for rows 2-5, regenerated to filter save games.
# Column F (Win) is left unchanged; column G (sum) = B + C + D + E.

$ws.Range("B2").Value = 0.3048080303191223
$ws.Range("C2").Value = 0.04240448674262143
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.164970295987679

$ws.Range("B3").Value = 1.459612070389937
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 15.68806981981553

$ws.Range("B4").Value = 0.04763786555579896
$ws.Range("C4").Value = 0.04240448674262143
$ws.Range("D4").Value = 0.1575252929769615
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0.7443468554461139

$ws.Range("B5").Value = 0.003994804209775715
$ws.Range("C5").Value = 0.04240448674262143
$ws.Range("D5").Value = 3.900430680208489
$ws.Range("E5").Value = 8.660232485948974
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 12.60706245710986
